# Timing issue fix - keywords, updated tc1,2 in ubc01
# Update the Cypher query text stored in column B for the Cases/Samples/Files tabs,
# reorder so Samples/Files text precedes Cases text in the shared-string table,
# drop the trailing Cohort line from the Cases query, and refresh the view selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$casesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)

MATCH (c)<--(diag:diagnosis)
	WHERE s.clinical_study_designation IN ['UBC01'] and demo.breed in ['Belgian Malinois', 'Labrador Retriever','West Highland White Terrier']and diag.disease_term in ['Bladder Cancer'] and diag.primary_disease_site in [ 'Bladder, Prostate', 'Bladder, Urethra','Bladder, Urethra, Prostate'] and diag.best_response in ['Partial Response', 'Progressive Disease']
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

$samplesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (samp:sample)-->(c)<--(diag:diagnosis) 
	WHERE s.clinical_study_designation IN ['UBC01'] and demo.breed in ['Belgian Malinois', 'Labrador Retriever','West Highland White Terrier']and diag.disease_term in ['Bladder Cancer'] and diag.primary_disease_site in [ 'Bladder, Prostate', 'Bladder, Urethra','Bladder, Urethra, Prostate'] and diag.best_response in ['Partial Response', 'Progressive Disease']
 WITH DISTINCT samp AS samp, c, demo, diag
RETURN  coalesce(samp.sample_id, '') AS `Sample ID`, 
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(demo.breed,'') AS Breed , 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(samp.sample_site, '') AS `Sample Site`,
        coalesce(samp.summarized_sample_type, '') AS `Sample Type`,
        coalesce(samp.specific_sample_pathology, '') AS `Pathology/Morphology`,
        coalesce(samp.tumor_grade, '') AS `Tumor Grade`,
        coalesce(samp.sample_chronology, '') AS `Sample Chronology`,
        coalesce(samp.percentage_tumor, '') AS `Percentage Tumor`,
        coalesce(samp.necropsy_sample, '') AS `Necropsy Sample`,
        coalesce(samp.sample_preservation, '') AS `Sample Preservation`
'@

$filesQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
 MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
	WHERE s.clinical_study_designation IN ['UBC01'] and demo.breed in ['Belgian Malinois', 'Labrador Retriever','West Highland White Terrier']and diag.disease_term in ['Bladder Cancer'] and diag.primary_disease_site in [ 'Bladder, Prostate', 'Bladder, Urethra','Bladder, Urethra, Prostate'] and diag.best_response in ['Partial Response', 'Progressive Disease']    
WITH DISTINCT f, parent, c, demo, diag, s
RETURN coalesce(f.file_name, '') AS `File Name`, 
        coalesce(f.file_type, '') AS `File Type`, 
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `File Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(demo.breed,'') AS Breed , 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

# Row 2 (CasesTab) now holds the shortened Cases query (no Cohort column)
$ws.Range("B2").Value = $casesQuery
$ws.Rows.Item(2).RowHeight = 319

# Row 3 (SamplesTab) holds the Samples query
$ws.Range("B3").Value = $samplesQuery

# Row 4 (FilesTab) holds the Files query
$ws.Range("B4").Value = $filesQuery

# Refresh the saved view/selection to match the re-saved workbook
$ws.Range("B2").Select() | Out-Null
